$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "phase" helper table in columns L/M, rows 4-8.
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 2
$ws.Range("M5").Formula = "=SQRT(L4^2+L5^2)"
$ws.Range("L6").Value = 3
$ws.Range("M6").Formula = "=SQRT(M5^2+L6^2)"
$ws.Range("L7").Value = 4
$ws.Range("M7").Formula = "=SQRT(M6^2+L7^2)"
$ws.Range("L8").Value = 5
$ws.Range("M8").Formula = "=SQRT(M7^2+L8^2)"
$ws.Range("N8").Formula = "=SQRT(L4^2+L5^2+L6^2+L8^2+L7^2)"

# C6 now left-aligned (new cellXfs entry applied).
$ws.Range("C6").HorizontalAlignment = -4131

# Move the active selection to just past the new data, matching the saved view.
[void]$ws.Range("N9").Select()
